$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.443991661071777
$ws.Range("B1").Value = 3.31175971031189
$ws.Range("C1").Value = 4.249300003051758
$ws.Range("D1").Value = 1.971782088279724
$ws.Range("E1").Value = 1.157291650772095
